# Reduce WHO Pertussis from 5-dose to 4-dose (3p+1) series per WHO guidance.
# - Rename the "5-dose series" sheet to "4-dose series"
# - Update the series-name cell text
# - Remove the now-unneeded "recurring dose" extra columns that were padding
#   the allowable-vaccine rows for a 5th dose
# - Remove the whole Dose 5 block (rows 42-50)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5-dose series")

# Update the series name text in B1 (shared string "WHO Pertussis 5-dose series"
# -> "WHO Pertussis 4-dose series (3p+1)")
$ws.Range("B1").Value = "WHO Pertussis 4-dose series (3p+1)"

# Trim the trailing "n/a" filler cells that followed each Dose-3 allowable
# vaccine row (these existed only to pad out to the old 5-dose width)
$ws.Range("G8:H8").Clear()
$ws.Range("G16:H16").Clear()
$ws.Range("J17:L17").Clear()
$ws.Range("G25:H25").Clear()
$ws.Range("J26:L26").Clear()
$ws.Range("G34:H34").Clear()
$ws.Range("J35:L35").Clear()

# Remove the entire Dose 5 block (rows 42-50): Series Dose/Age/Allowable
# Vaccine rows for the school-entry booster WHO no longer recommends.
$ws.Range("A42:L50").Clear()

# Rename the sheet itself last, after all range lookups by the old name.
$ws.Name = "4-dose series"
